$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeouts) values replacing the old "Strike#" derived values in column G.
$newK = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 1
    6  = 2
    7  = 0
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 2
    14 = 0
    15 = 1
    16 = 2
    17 = 1
    18 = 1
    19 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
